$d = $word.ActiveDocument

$d.Content.Find.Execute("78×62=4836", $true, $false, $false, $false, $false, $true, 1, $false, "13×75=975", 2)
$d.Content.Find.Execute("66×30=1980", $true, $false, $false, $false, $false, $true, 1, $false, "96×85=8160", 2)
$d.Content.Find.Execute("57×23=1311", $true, $false, $false, $false, $false, $true, 1, $false, "59×70=4130", 2)
$d.Content.Find.Execute("72×15=1080", $true, $false, $false, $false, $false, $true, 1, $false, "57×62=3534", 2)
$d.Content.Find.Execute("57×26=1482", $true, $false, $false, $false, $false, $true, 1, $false, "99×29=2871", 2)
$d.Content.Find.Execute("60×14=840", $true, $false, $false, $false, $false, $true, 1, $false, "81×37=2997", 2)
$d.Content.Find.Execute("23×84=1932", $true, $false, $false, $false, $false, $true, 1, $false, "93×52=4836", 2)
$d.Content.Find.Execute("56×22=1232", $true, $false, $false, $false, $false, $true, 1, $false, "37×19=703", 2)
$d.Content.Find.Execute("39×95=3705", $true, $false, $false, $false, $false, $true, 1, $false, "95×96=9120", 2)
$d.Content.Find.Execute("77×72=5544", $true, $false, $false, $false, $false, $true, 1, $false, "48×72=3456", 2)
$d.Content.Find.Execute("33×87=2871", $true, $false, $false, $false, $false, $true, 1, $false, "29×38=1102", 2)
$d.Content.Find.Execute("42×14=588", $true, $false, $false, $false, $false, $true, 1, $false, "33×88=2904", 2)
$d.Content.Find.Execute("40×77=3080", $true, $false, $false, $false, $false, $true, 1, $false, "90×45=4050", 2)
$d.Content.Find.Execute("78×35=2730", $true, $false, $false, $false, $false, $true, 1, $false, "56×33=1848", 2)
$d.Content.Find.Execute("83×38=3154", $true, $false, $false, $false, $false, $true, 1, $false, "40×28=1120", 2)
$d.Content.Find.Execute("43×63=2709", $true, $false, $false, $false, $false, $true, 1, $false, "14×13=182", 2)
$d.Content.Find.Execute("28×55=1540", $true, $false, $false, $false, $false, $true, 1, $false, "74×90=6660", 2)
$d.Content.Find.Execute("56×65=3640", $true, $false, $false, $false, $false, $true, 1, $false, "14×94=1316", 2)
$d.Content.Find.Execute("93×73=6789", $true, $false, $false, $false, $false, $true, 1, $false, "73×33=2409", 2)
$d.Content.Find.Execute("58×82=4756", $true, $false, $false, $false, $false, $true, 1, $false, "52×85=4420", 2)
$d.Content.Find.Execute("39×54=2106", $true, $false, $false, $false, $false, $true, 1, $false, "97×35=3395", 2)
$d.Content.Find.Execute("92×91=8372", $true, $false, $false, $false, $false, $true, 1, $false, "30×60=1800", 2)
$d.Content.Find.Execute("52×26=1352", $true, $false, $false, $false, $false, $true, 1, $false, "87×91=7917", 2)
$d.Content.Find.Execute("12×77=924", $true, $false, $false, $false, $false, $true, 1, $false, "33×37=1221", 2)
$d.Content.Find.Execute("68×70=4760", $true, $false, $false, $false, $false, $true, 1, $false, "66×67=4422", 2)
